# it279-program5-pseudocode.docx edit script
# Implements (per commit diff):
#  1. Insert a new bold/52pt "dijkstra's algorithm is on next page" paragraph
#     right before the page-break paragraph that precedes the Dijkstra section.
#  2. Rename the AdjListVertex struct to Edge and rename/retype its fields
#     (destination/edgeCost -> toIndex/cost) in the struct definition.
#  3. Rename "for each AdjListVertex curAdjVertex adjacent to startingVertex"
#     to use "Edge" instead of "AdjListVertex".
#  4. Merge the split "starting" / "Vertex." / "totalDistance" runs (and the
#     following " " / "+ " runs) into single runs.
#  5. Rename "for each AdjListVertex(split as A/djListVertex) curAdjVertex
#     adjacent to curQueueVertex.toVertex" to use a single "Edge" run.

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParaIndex {
    param($doc, [string]$substr)
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Contains($substr)) {
            return $i
        }
    }
    return -1
}

function Find-ParaIndexWithBreak {
    param($doc)
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        if ($doc.Paragraphs($i).Range.WordOpenXML -match 'w:br') {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Insert the new "dijkstra's algorithm is on next page" paragraph just
#    before the paragraph holding the manual page break.
# ---------------------------------------------------------------------------
$breakIdx = Find-ParaIndexWithBreak $d
$d.Paragraphs($breakIdx).Range.InsertParagraphBefore()

$newParaIdx = $breakIdx
$newParaXml = @"
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="52"/>
      <w:szCs w:val="52"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="52"/>
      <w:szCs w:val="52"/>
    </w:rPr>
    <w:t>dijkstra&#8217;s</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="52"/>
      <w:szCs w:val="52"/>
    </w:rPr>
    <w:t xml:space="preserve"> algorithm is on next </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="52"/>
      <w:szCs w:val="52"/>
    </w:rPr>
    <w:t>page</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
[void]$d.Paragraphs($newParaIdx).Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 2. "AdjListVertex: int destination, int edgeCost" -> "Edge: int toIndex, int cost"
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "AdjListVertex: int destination"
$xml = @"
<w:p $wns>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:t>Edge</w:t></w:r>
  <w:r><w:t xml:space="preserve">: </w:t></w:r>
  <w:r><w:t xml:space="preserve">int </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>toIndex</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, int </w:t></w:r>
  <w:r><w:t>c</w:t></w:r>
  <w:r><w:t>ost</w:t></w:r>
</w:p>
"@
[void]$d.Paragraphs($idx).Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3. "for each AdjListVertex curAdjVertex adjacent to startingVertex" -> Edge
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "for each AdjListVertex curAdjVertex adjacent to startingVertex"
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:vertAlign w:val="subscript"/>
    </w:rPr>
  </w:pPr>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/><w:t xml:space="preserve">for each </w:t></w:r>
  <w:r><w:t>Edge</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>curAdjVertex</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> adjacent to </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>startingVertex</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
[void]$d.Paragraphs($idx).Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4. Merge "starting"+"Vertex."+"totalDistance" and " "+"+ " runs together
#    in the "totalCost as  startingVertex.totalDistance + curAdjacentVertex.edgeCost"
#    paragraph.
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "totalCost as  startingVertex.totalDistance"
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:ind w:left="2160" w:firstLine="720"/>
    <w:rPr>
      <w:sz w:val="16"/>
      <w:szCs w:val="16"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>totalCost</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> as  </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>startingVertex.totalDistance</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> + </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>curAdjacentVertex.edgeCost</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
[void]$d.Paragraphs($idx).Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5. "for each AdjListVertex curAdjVertex adjacent to curQueueVertex.toVertex" -> Edge
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "for each AdjListVertex curAdjVertex adjacent to curQueueVertex.toVertex"
$xml = @"
<w:p $wns>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:t xml:space="preserve">for each </w:t></w:r>
  <w:r><w:t>Edge</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>cur</w:t></w:r>
  <w:r><w:t>AdjVertex</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> adjacent to </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>curQueueVertex</w:t></w:r>
  <w:r><w:t>.toVertex</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
[void]$d.Paragraphs($idx).Range.InsertXML($xml)

Write-Output "done"
